# [draft] Adding exiftool + GTN "new version check" -- untested.
#
# Appends 4 new rows (40-43) to the "messageBox" sheet, each a pair of
# (message-key, message-text) cells in columns A/C, mirroring the existing
# rows on the sheet. This introduces 8 new shared strings:
#   - mbx_Helper_WarningExifToolVerAPIResponse
#   - "ExifToolVerAPIResponse API Returned the following response: "
#   - mbx_frm_mainApp_InfoNewExifToolVersionExists
#   - "There is a new exifTool version available for download. ..."
#   - mbx_Helper_WarningGTNVerAPIResponse
#   - "GitHub API Returned the following response: "
#   - mbx_frm_mainApp_InfoNewGTNVersionExists
#   - "There is a new GeoTagNinja version available for download. ..."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40: exiftool version-check API-response warning -------------------
$ws.Range("A40").Value = "mbx_Helper_WarningExifToolVerAPIResponse"
$ws.Range("C40").Value = "ExifToolVerAPIResponse API Returned the following response: "

# --- Row 41: exiftool new-version-available info message -------------------
$ws.Range("A41").Value = "mbx_frm_mainApp_InfoNewExifToolVersionExists"
$ws.Range("C41").Value = "There is a new exifTool version available for download. `nIf you click Yes, your default browser will open to exiftool.org where you can download manually. If you click No, this message will close (but will show again next time.)`nIf you download, extract the zip file, rename exiftool(-k).exe to exiftool.exe and replace the current file in your GeoTagNinja installation folder (most likely Program Files.)`nNew version: "

# --- Row 42: GeoTagNinja version-check API-response warning ----------------
$ws.Range("A42").Value = "mbx_Helper_WarningGTNVerAPIResponse"
$ws.Range("C42").Value = "GitHub API Returned the following response: "

# --- Row 43: GeoTagNinja new-version-available info message ----------------
$ws.Range("A43").Value = "mbx_frm_mainApp_InfoNewGTNVersionExists"
$ws.Range("C43").Value = "There is a new GeoTagNinja version available for download. `nIf you click Yes, your default browser will open directly to the release's installer link on GitHub, which you can then install manually. If you click No, this message will close (but will show again next time.)`nNew version: "

# The two multi-line messages wrap onto several lines in column C (which is
# styled with wrapText); match the row heights Excel computed for them.
$ws.Rows.Item(41).RowHeight = 116.6
$ws.Rows.Item(43).RowHeight = 72.9

# Leave the sheet's selection on the newly-added last cell, as in the diff.
$null = $ws.Range("C43").Select()
